$d = $word.ActiveDocument

# Locate the paragraph containing the "LOB1024: Mecânica (Requisito fraco)" text
# and the paragraph containing the trailing empty paragraph right before the
# page-break paragraph. We need to delete everything in between (inclusive of
# the blank paragraph, the "Ver no Jupiter..." paragraph and the "© 2020..."
# paragraph), but keep the "LOB1024" paragraph itself and the paragraph mark
# that follows it as-is.

$n = $d.Paragraphs.Count
$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*LOB1024: Mec*(Requisito fraco)*") {
        $startIndex = $i
    }
    if ($t -like "*2020 . Contact: luizeleno@usp.br*") {
        $endIndex = $i
    }
}

$deleteStart = $d.Paragraphs.Item($startIndex).Range.End
$deleteEnd = $d.Paragraphs.Item($endIndex).Range.End

$r = $d.Range($deleteStart, $deleteEnd)
$r.Delete()
